# =====================================================================
# edit.ps1 -- Gungnir_Profits workbook refresh (scheduled runner update)
#
# Refreshes the derived market-price / leve-profit columns
# (currentAveragePrice, currentAveragePriceNQ, currentAveragePriceHQ,
#  LevePriceNQ, LevePriceHQ, LeveProfitNQ, LeveProfitHQ -- columns H:N)
# for the affected leve rows on each job sheet, matching a refreshed
# market-data pull. Values only; no structural changes.
# =====================================================================

$wb = $excel.ActiveWorkbook

function Set-LeveRow {
    param($ws, $Row, $Values)
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value2 = $Values[$col]
    }
}

# --- ALC ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")
Set-LeveRow $ws 64 @{ "H" = 3256.6667; "I" = 3328.5715; "J" = 3193.75; "K" = 3328.5715; "L" = 3193.75; "M" = -3080.5715; "N" = -3689.75 }
Set-LeveRow $ws 67 @{ "H" = 3256.6667; "I" = 3328.5715; "J" = 3193.75; "K" = 3328.5715; "L" = 3193.75; "M" = -2470.5715; "N" = -4909.75 }
Set-LeveRow $ws 129 @{ "H" = 1043.5636; "I" = 1215.4286; "J" = 1018.5; "K" = 3646.2858; "L" = 3055.5; "M" = 1353.7142; "N" = -13055.5 }
Set-LeveRow $ws 138 @{ "H" = 2258.8235; "I" = 1511.6; "K" = 4534.799999999999; "M" = 605.2000000000007 }

# --- ARM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")
Set-LeveRow $ws 8 @{ "H" = 1599.3334; "I" = 290; "J" = 1861.2; "K" = 290; "L" = 1861.2; "M" = -146; "N" = -2149.2 }
Set-LeveRow $ws 21 @{ "H" = 2966.6667; "J" = 4500; "L" = 4500; "N" = -5248 }
Set-LeveRow $ws 26 @{ "H" = 840.7778; "I" = 652.4286; "J" = 1500; "K" = 652.4286; "L" = 1500; "M" = -322.4286; "N" = -2160 }
Set-LeveRow $ws 34 @{ "H" = 28846.666; "J" = 37693.332; "L" = 37693.332; "N" = -38235.332 }

# --- BSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")
Set-LeveRow $ws 8 @{ "H" = 280.8; "I" = 101; "J" = 1000; "K" = 101; "L" = 1000; "M" = 39; "N" = -1280 }
Set-LeveRow $ws 14 @{ "H" = 270; "I" = 320; "J" = 220; "K" = 320; "L" = 220; "M" = -148; "N" = -564 }
Set-LeveRow $ws 105 @{ "H" = 100002140; "I" = 2064.1428; "J" = 333335680; "K" = 2064.1428; "L" = 333335680; "M" = -317.1428000000001; "N" = -333339174 }
Set-LeveRow $ws 107 @{ "H" = 25000738; "I" = 25000738; "K" = 25000738; "M" = -24998818 }

# --- CRP ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")
Set-LeveRow $ws 32 @{ "H" = 6974; "I" = 900; "K" = 900; "M" = -584 }
Set-LeveRow $ws 36 @{ "H" = 4764; "I" = 4449.3335; "J" = 5000; "K" = 4449.3335; "L" = 5000; "M" = -4061.3335; "N" = -5776 }
Set-LeveRow $ws 38 @{ "H" = 8520.791999999999; "I" = 3136.2727; "J" = 13076.923; "K" = 3136.2727; "L" = 13076.923; "M" = -2759.2727; "N" = -13830.923 }
Set-LeveRow $ws 40 @{ "H" = 4764; "I" = 4449.3335; "J" = 5000; "K" = 4449.3335; "L" = 5000; "M" = -4289.3335; "N" = -5320 }
Set-LeveRow $ws 46 @{ "H" = 8520.791999999999; "I" = 3136.2727; "J" = 13076.923; "K" = 3136.2727; "L" = 13076.923; "M" = -2925.2727; "N" = -13498.923 }

# --- CUL ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")
Set-LeveRow $ws 23 @{ "H" = 109.17647; "J" = 109.17647; "L" = 327.52941; "N" = -797.52941 }
Set-LeveRow $ws 115 @{ "H" = 4222.615; "I" = 475.5; "J" = 5888; "K" = 1426.5; "L" = 17664; "M" = -251.5; "N" = -20014 }
Set-LeveRow $ws 137 @{ "H" = 11906069; "I" = 26316492; "K" = 78949476; "M" = -78944376 }

# --- GSM ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")
Set-LeveRow $ws 14 @{ "H" = 213.33333; "I" = 213.33333; "K" = 213.33333; "M" = -45.33332999999999 }
Set-LeveRow $ws 102 @{ "H" = 1347.1515; "I" = 746.5185; "J" = 4050; "K" = 746.5185; "L" = 4050; "M" = 875.4815; "N" = -7294 }
Set-LeveRow $ws 113 @{ "H" = 1301; "I" = 1151.25; "J" = 1900; "K" = 1151.25; "L" = 1900; "M" = 1018.75; "N" = -6240 }

# --- LTW ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")
Set-LeveRow $ws 18 @{ "H" = 1099; "I" = 1099; "K" = 1099; "M" = -927 }
Set-LeveRow $ws 20 @{ "H" = 8333.333000000001; "J" = 8333.333000000001; "L" = 8333.333000000001; "N" = -8785.333000000001 }
Set-LeveRow $ws 21 @{ "H" = 1000; "I" = 1000; "K" = 1000; "M" = -826 }
Set-LeveRow $ws 22 @{ "H" = 1057039; "J" = 2151.5715; "L" = 2151.5715; "N" = -2741.5715 }
Set-LeveRow $ws 26 @{ "H" = 1000; "I" = 1000; "K" = 1000; "M" = -705 }
Set-LeveRow $ws 27 @{ "H" = 1057039; "J" = 2151.5715; "L" = 2151.5715; "N" = -2365.5715 }
Set-LeveRow $ws 47 @{ "H" = 12000; "J" = 12000; "L" = 12000; "N" = -12980 }
Set-LeveRow $ws 52 @{ "H" = 12000; "J" = 12000; "L" = 12000; "N" = -12466 }

# --- WVR ---------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")
Set-LeveRow $ws 107 @{ "H" = 499.875; "I" = 499.8; "J" = 500; "K" = 1499.4; "L" = 1500; "M" = 420.5999999999999; "N" = -5340 }
Set-LeveRow $ws 122 @{ "H" = 38241.355; "I" = 51537.9; "J" = 5000; "K" = 154613.7; "L" = 15000; "M" = -152163.7; "N" = -19900 }

